$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty price cells D2 and D3
$ws.Range("D2").Value = 178
$ws.Range("D3").Value = 178

# Scroll the sheet view back to the top-left and move the selection to D4
# (undoes the previous scrolled-down view / D34 selection state)
[void]$ws.Range("A1").Select()
[void]$ws.Range("D4").Select()
